$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A12").Value = 0
$ws.Range("A49").Value = 1
$ws.Range("A52").Value = 1
$ws.Range("A53").Value = 1
$ws.Range("A58").Value = 1
$ws.Range("A80").Value = 0
$ws.Range("A156").Value = 1
$ws.Range("A258").Value = 0
$ws.Range("A278").Value = 1
$ws.Range("A351").Value = 0
$ws.Range("A376").Value = 0
$ws.Range("A481").Value = 1
$ws.Range("A503").Value = 1
